$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the F:G columns which are being removed entirely
$ws.Range("F1:G3").Clear()

# Update header row (B1:E1) with the new, reordered label set
$ws.Range("B1").Value = "C_B"
$ws.Range("C1").Value = "C_LF"
$ws.Range("D1").Value = "FFR_B"
$ws.Range("E1").Value = "FFR_LF"

# Update params row (row 2) values
$ws.Range("B2").Value = 0.87504171261757
$ws.Range("C2").Value = 0.00337956661413657
$ws.Range("D2").Value = -18.0515008959796
$ws.Range("E2").Value = 0.4919540470472872

# Update pvalue row (row 3) values
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = [double]"5.788347579027686E-11"
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
